# Correção nos subtítulos da planilha de demonstrativo financeiro.
# Replace hyphen ("-") with en-dash ("–") and fix stray whitespace in a
# set of numbered subtitle/caption cells on the "PTRF Básico" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PTRF Básico")

$ws.Range("D14").Value = "10 – Rendimento"
$ws.Range("A19").Value = "16 – Item"
$ws.Range("B20").Value = "17 – Razão social"
$ws.Range("F20").Value = "21 – Data"
$ws.Range("G20").Value = "22 – Especificação do material ou serviço"
$ws.Range("H20").Value = "23 - Tipo da despesa "
$ws.Range("I20").Value = "24 -  Tipo de transação"
$ws.Range("A24").Value = "27 – Item"
$ws.Range("B25").Value = "28 – Razão social"
$ws.Range("F25").Value = "32 – Data"
$ws.Range("G25").Value = "33 – Especificação do material ou serviço"
$ws.Range("H25").Value = "34 - Tipo da despesa "
$ws.Range("I25").Value = "35 -  Tipo de transação"
$ws.Range("A29").Value = "38 – Item"
$ws.Range("F29").Value = "40 – Data"
$ws.Range("H29").Value = "41 - Valor"

# Row-height touch-ups that accompanied the subtitle fixes (rows shrank from
# the old 25.5pt autofit height down to the sheet's standard 21pt row).
$ws.Rows.Item(20).RowHeight = 21
$ws.Rows.Item(25).RowHeight = 21
$ws.Rows.Item(29).RowHeight = 21
